$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.400.72"
$ws.Range("E2").Value = "  -2.57%  "

$ws.Range("D3").Value = "3.177.77"
$ws.Range("E3").Value = "  -4.22%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'570.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.59%  "

$ws.Range("D6").Value = "'168.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.29%  "

$ws.Range("D7").Value = "'0.608"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -6.90%  "

$ws.Range("E8").Value = "  -0.15%  "

$ws.Range("D9").Value = "3.187.39"
$ws.Range("E9").Value = "  -3.91%  "

$ws.Range("E10").Value = "  -4.14%  "

$ws.Range("E11").Value = "  -0.22%  "

$ws.Range("D12").Value = "'0.389"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.17%  "

$ws.Range("D13").Value = "3.733.36"
$ws.Range("E13").Value = "  -4.06%  "

$ws.Range("E14").Value = "  -1.45%  "

$ws.Range("D15").Value = "64.477.17"
$ws.Range("E15").Value = "  -2.56%  "

$ws.Range("D16").Value = "'25.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.39%  "

$ws.Range("E17").Value = "  -3.82%  "

$ws.Range("D18").Value = "3.191.97"
$ws.Range("E18").Value = "  -1.94%  "

$ws.Range("D19").Value = "'420.20"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.11%  "

$ws.Range("D20").Value = "'12.99"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.06%  "

$ws.Range("E21").Value = "  -3.35%  "

$ws.Range("D22").Value = "'7.14"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.16%  "

$ws.Range("E23").Value = "  -0.17%  "

$ws.Range("B24").Value = "LEO"
$ws.Range("C24").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D24").Value = "'5.67"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.03%  "

$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "'70.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.13%  "

$ws.Range("D26").Value = "'0.205"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.44%  "

$ws.Range("D27").Value = "'0.496"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.07%  "

$ws.Range("E28").Value = "  -7.31%  "

$ws.Range("D29").Value = "'8.72"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.95%  "

$ws.Range("D30").Value = "'0.998"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.32%  "

$ws.Range("D31").Value = "'1.83"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.52%  "

$ws.Range("D32").Value = "'21.78"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.71%  "

$ws.Range("D33").Value = "'0.999"
$ws.Range("D33").Style = "Normal"

$ws.Range("D34").Value = "'5.04"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.55%  "

$ws.Range("E35").Value = "  -3.15%  "

$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").Value = "'157.13"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.05%  "

$ws.Range("B37").Value = "Fetch.AI"
$ws.Range("C37").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D37").Value = "'1.12"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.44%  "

$ws.Range("E38").Value = "  -4.86%  "

$ws.Range("E39").Value = "  -5.27%  "

$ws.Range("D40").Value = "2.690.84"
$ws.Range("E40").Value = "  -6.54%  "

$ws.Range("D41").Value = "'4.23"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.74%  "

$ws.Range("E42").Value = "  -8.12%  "

$ws.Range("D43").Value = "'39.28"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.37%  "

$ws.Range("E44").Value = "  -5.68%  "

$ws.Range("E45").Value = "  -5.66%  "

$ws.Range("D46").Value = "'5.59"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.37%  "

$ws.Range("E47").Value = "  -2.94%  "

$ws.Range("D48").Value = "'291.70"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -7.05%  "

$ws.Range("E49").Value = "  -7.69%  "

$ws.Range("D50").Value = "'0.0995"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.15%  "

$ws.Range("E51").Value = "  -0.19%  "
